$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume figures for this data refresh.
# "Price" values are plain text in this sheet (not real numbers - note the
# thousands-dot formatting like "45.057.60"), but many of the new quotes are
# numeric-looking ("302.73", "0.800", ...). Excel auto-detects those as numbers
# on a plain .Value assignment (and can even drop significant trailing zeros,
# e.g. "238.30" -> 238.3), so every such cell is pinned to Text format right
# before the write to keep it a string and preserve the exact display text.

$ws.Range("D2").Value = '45.057.60'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '2.265.17'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("E4").Value = '  -0.72%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.73'
$ws.Range("E5").Value = '  -1.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.27'
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("E7").Value = '  -1.61%  '
$ws.Range("E8").Value = '  -0.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.511'
$ws.Range("E9").Value = '  -2.16%  '
$ws.Range("E10").Value = '  -3.33%  '
$ws.Range("E11").Value = '  -2.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.21'
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("E13").Value = '  -1.07%  '
$ws.Range("D14").Value = '2.610.21'
$ws.Range("E14").Value = '  +0.68%  '
$ws.Range("D15").Value = '2.271.19'
$ws.Range("E15").Value = '  -1.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.62'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.800'
$ws.Range("E17").Value = '  -5.13%  '
$ws.Range("D18").Value = '44.909.67'
$ws.Range("E18").Value = '  +0.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.94'
$ws.Range("E19").Value = '  +7.62%  '
$ws.Range("D20").Value = '0.0₃0924'
$ws.Range("E20").Value = '  -2.86%  '
$ws.Range("E21").Value = '  -3.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.66'
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.30'
$ws.Range("E23").Value = '  -0.87%  '
$ws.Range("E24").Value = '  -2.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("E25").Value = '  -0.49%  '
$ws.Range("E26").Value = '  -5.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '41.53'
$ws.Range("E27").Value = '  +10.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.26'
$ws.Range("E28").Value = '  -0.80%  '
$ws.Range("E29").Value = '  -2.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.57'
$ws.Range("E30").Value = '  -2.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '153.10'
$ws.Range("E31").Value = '  +1.56%  '
$ws.Range("E32").Value = '  -8.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0789'
$ws.Range("E33").Value = '  -1.68%  '
$ws.Range("E34").Value = '  -2.64%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.96'
$ws.Range("E35").Value = '  -3.15%  '
$ws.Range("E36").Value = '  -1.75%  '
$ws.Range("E37").Value = '  -3.56%  '
$ws.Range("E38").Value = '  -5.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.05'
$ws.Range("E39").Value = '  +6.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0315'
$ws.Range("E40").Value = '  +4.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.26'
$ws.Range("E41").Value = '  -3.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.65'
$ws.Range("E42").Value = '  -8.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  -0.86%  '
$ws.Range("E44").Value = '  +12.25%  '
$ws.Range("D45").Value = '1.744.10'
$ws.Range("E45").Value = '  -5.08%  '
$ws.Range("E46").Value = '  +3.94%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '71.06'
$ws.Range("E47").Value = '  +2.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '76.07'
$ws.Range("E48").Value = '  -5.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '96.37'
$ws.Range("E49").Value = '  -3.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.98'
$ws.Range("E50").Value = '  -1.50%  '
$ws.Range("E51").Value = '  -4.70%  '
